$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Flt3l", "Flt3", "Inflammatory-Mac", 2, 1, 7.6494085, 15.298817, 0.2611121155001731, 0.2017467183721511, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 10.85758827351983, 65.145529641119, 0.1066003261861801, 0.08236410610157406),
    @("ECs", "Flt3l", "Flt3", "Neutrophils", 2, 1, 7.6494085, 15.298817, 0.2611121155001731, 0.2017467183721511, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 11.16309289999567, 66.97855739997399, 0.1095997853674734, 0.08468162034467218),
    @("ECs", "Flt3l", "Flt3", "Resolving-Mac", 2, 1, 7.6494085, 15.298817, 0.2611121155001731, 0.2017467183721511, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 4.574432976296333, 27.446597857778, 0.04491200394651962, 0.0347009919259048),
    @("FAPs", "Flt3l", "Flt3", "Inflammatory-Mac", 3, 1, 3.299743666666667, 9.899231, 0.1126365586844566, 0.1305419476981696, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 4.683663859868555, 42.152974738817, 0.04598443803824161, 0.05329440259387318),
    @("FAPs", "Flt3l", "Flt3", "Neutrophils", 3, 1, 3.299743666666667, 9.899231, 0.1126365586844566, 0.1305419476981696, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 4.815450122342444, 43.339051101082, 0.04727832192964283, 0.05479397009887821),
    @("FAPs", "Flt3l", "Flt3", "Resolving-Mac", 3, 1, 3.299743666666667, 9.899231, 0.1126365586844566, 0.1305419476981696, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 1.973284109761556, 17.759556987854, 0.0193737987165722, 0.02245357500541817),
    @("Inflammatory-Mac", "Flt3l", "Flt3", "Inflammatory-Mac", 3, 1, 6.035325333333333, 18.105976, 0.2060154802189547, 0.2387649477031411, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 8.56655486055911, 77.09899374503199, 0.08410684946071968, 0.09747698324233522),
    @("Inflammatory-Mac", "Flt3l", "Flt3", "Neutrophils", 3, 1, 6.035325333333333, 18.105976, 0.2060154802189547, 0.2387649477031411, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 8.807595695496888, 79.26836125947199, 0.08647340002252565, 0.1002197350031539),
    @("Inflammatory-Mac", "Flt3l", "Flt3", "Resolving-Mac", 3, 1, 6.035325333333333, 18.105976, 0.2060154802189547, 0.2387649477031411, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 3.609192949687111, 32.482736547184, 0.03543523073570937, 0.04106822945765193),
    @("MuSCs", "Flt3l", "Flt3", "Inflammatory-Mac", 2, 1, 4.4052735, 8.810547, 0.150373755427214, 0.1161853850734734, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 6.2528554848715, 37.517132909229, 0.06139083721824182, 0.0474332641485224),
    @("MuSCs", "Flt3l", "Flt3", "Neutrophils", 2, 1, 4.4052735, 8.810547, 0.150373755427214, 0.1161853850734734, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 6.428794766338999, 38.572768598034, 0.06311821758310047, 0.04876791428271156),
    @("MuSCs", "Flt3l", "Flt3", "Resolving-Mac", 2, 1, 4.4052735, 8.810547, 0.150373755427214, 0.1161853850734734, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 2.634403479433, 15.806420876598, 0.0258647006258717, 0.01998420664223938),
    @("Neutrophils", "Flt3l", "Flt3", "Inflammatory-Mac", 3, 1, 4.902894666666666, 14.708684, 0.1673600250905478, 0.1939645874954229, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 6.959180129954221, 62.63262116958798, 0.06832556670534061, 0.07918701227621225),
    @("Neutrophils", "Flt3l", "Flt3", "Neutrophils", 3, 1, 4.902894666666666, 14.708684, 0.1673600250905478, 0.1939645874954229, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 7.154993571449777, 64.39494214304798, 0.07024807253345099, 0.08141513126523142),
    @("Neutrophils", "Flt3l", "Flt3", "Resolving-Mac", 3, 1, 4.902894666666666, 14.708684, 0.1673600250905478, 0.1939645874954229, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 2.931986576806222, 26.38787919125599, 0.02878638585175616, 0.03336244395397926),
    @("Resolving-Mac", "Flt3l", "Flt3", "Inflammatory-Mac", 3, 1, 3.002848666666667, 9.008546000000001, 0.1025020650786538, 0.1187964136576422, 3, 1, 1.419402333333333, 4.258207, 0.4082549979803961, 0.4082549979803961, 4.262250404113556, 38.360253637022, 0.04184698037167223, 0.048499229617879),
    @("Resolving-Mac", "Flt3l", "Flt3", "Neutrophils", 3, 1, 3.002848666666667, 9.008546000000001, 0.1025020650786538, 0.1187964136576422, 3, 1, 1.459340666666667, 4.378022, 0.4197422442751443, 0.4197422442751443, 4.382179175112444, 39.439612576012, 0.04302444683895105, 0.04986387328049713),
    @("Resolving-Mac", "Flt3l", "Flt3", "Resolving-Mac", 3, 1, 3.002848666666667, 9.008546000000001, 0.1025020650786538, 0.1187964136576422, 3, 1, 0.5980113333333333, 1.794034, 0.1720027577444596, 0.1720027577444596, 1.795737534951556, 16.161637814564, 0.01763063786803052, 0.02043331075926603)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $c = 1 + $j
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
